# Generate Report for Handback
# - Marks Overview + per-locale Status as "Handed back: in sync with en-US"
# - Fills in "Latest Target File" (hyperlinked source doc), "Latest Handback File"
#   (generated xliff) and "Latest Handback DateTime" for each locale sheet's rows
# - Widens the columns that now hold the longer generated values

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

$mdUrl1820 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/44a5bb7bb5c03a042accccf1db8113adeeb2abdc/e2e/1820e3cb-b105-4c99-968c-e0a70946fd4d.md"
$mdUrl798d = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/44a5bb7bb5c03a042accccf1db8113adeeb2abdc/e2e/798d0941-2185-41e8-8db3-66ec1f5541b2.md"

$mdName1820 = "1820e3cb-b105-4c99-968c-e0a70946fd4d.md"
$mdName798d = "798d0941-2185-41e8-8db3-66ec1f5541b2.md"

# ---------------------------------------------------------------------------
# Overview sheet: update the per-locale handoff/handback status banner
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$overview.Range("E3").Value = $statusText
$overview.Range("F3").Value = $statusText

$overview.Columns.Item(5).ColumnWidth = 29.1
$overview.Columns.Item(6).ColumnWidth = 29.1

# ---------------------------------------------------------------------------
# zh-cn sheet: handback report for both source files
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = $statusText
$zhcn.Range("C3").Value = $statusText

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $mdUrl1820, "", "", $mdName1820)
$zhcn.Range("J2").Value = "1820e3cb-b105-4c99-968c-e0a70946fd4d.fd6b8ee2aa9675fe6f4677a9b24b0e03fe6a6c34.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-10-21 04:54:20"

$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $mdUrl798d, "", "", $mdName798d)
$zhcn.Range("J3").Value = "798d0941-2185-41e8-8db3-66ec1f5541b2.c7934fa8cb32ca8cd51ee993a4293c34793d8b39.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-10-21 04:54:20"

$zhcn.Columns.Item(3).ColumnWidth = 29.1
$zhcn.Columns.Item(9).ColumnWidth = 39.15
$zhcn.Columns.Item(10).ColumnWidth = 39.15

# ---------------------------------------------------------------------------
# de-de sheet: handback report for both source files
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = $statusText
$dede.Range("C3").Value = $statusText

$dede.Hyperlinks.Add($dede.Range("I2"), $mdUrl1820, "", "", $mdName1820)
$dede.Range("J2").Value = "1820e3cb-b105-4c99-968c-e0a70946fd4d.fd6b8ee2aa9675fe6f4677a9b24b0e03fe6a6c34.de-de.xlf"
$dede.Range("K2").Value = "2016-10-21 04:54:38"

$dede.Hyperlinks.Add($dede.Range("I3"), $mdUrl798d, "", "", $mdName798d)
$dede.Range("J3").Value = "798d0941-2185-41e8-8db3-66ec1f5541b2.c7934fa8cb32ca8cd51ee993a4293c34793d8b39.de-de.xlf"
$dede.Range("K3").Value = "2016-10-21 04:54:38"

$dede.Columns.Item(3).ColumnWidth = 29.1
$dede.Columns.Item(9).ColumnWidth = 39.15
$dede.Columns.Item(10).ColumnWidth = 39.15
